# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# These two sheets contain duplicated data, so the same row/value changes
# are applied to both.

$wb = $excel.ActiveWorkbook

$changes = @{
    3  = 5601
    4  = 38
    5  = 680
    6  = 674
    11 = 1588
    12 = 5435
    13 = 464
    14 = 290
    15 = 247
    19 = 4577
    21 = 1213
    22 = 123
    23 = 83
    26 = 205
    28 = 155
    29 = 85
    35 = 33
    37 = 45
    38 = 51
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Range("F$row").Value = $changes[$row]
    }
}
